$wb = $excel.ActiveWorkbook

$simNames = $wb.Worksheets.Item("Simulation Names")
$simNames.Range("D1").Value = "Lincoln2012NitNilIrrFull"
$simNames.Range("D2").Value = "Lincoln2012NitNilIrrNil"
$simNames.Range("D3").Value = "Lincoln2012NitLowIrrFull"
$simNames.Range("D4").Value = "Lincoln2012NitLowIrrNil"
$simNames.Range("D5").Value = "Lincoln2012NitMedIrrFull"
$simNames.Range("D6").Value = "Lincoln2012NitMedIrrNil"

$soilWater = $wb.Worksheets.Item("SoilWater")
$soilWater.Columns.Item(1).AutoFit() | Out-Null

$soilMoisturePerc = $wb.Worksheets.Item("soilMoisturePerc")
$soilMoistureMM = $wb.Worksheets.Item("SoilMoistureMM")

$soilWater.Activate()
$soilWater.Range("A2").Select() | Out-Null

$simNames.Activate()
$simNames.Range("E6:F6").Select() | Out-Null

$soilMoistureMM.Activate()
$soilMoistureMM.Range("D4").Select() | Out-Null

$soilMoisturePerc.Activate()
